$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.868.60"
$ws.Range("E2").Value = "  -0.20%  "
$ws.Range("D3").Value = "3.142.44"
$ws.Range("E3").Value = "  +0.39%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "591.97"
$ws.Range("E5").Value = "  +0.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.48"
$ws.Range("E6").Value = "  -1.25%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "3.133.49"
$ws.Range("E8").Value = "  +0.28%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.529"
$ws.Range("E9").Value = "  -0.68%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.160"
$ws.Range("E10").Value = "  -1.19%  "
$ws.Range("E11").Value = "  +1.89%  "
$ws.Range("E12").Value = "  -1.58%  "
$ws.Range("E13").Value = "  -3.14%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.14"
$ws.Range("E14").Value = "  -0.78%  "
$ws.Range("D15").Value = "3.663.20"
$ws.Range("E15").Value = "  +0.34%  "
$ws.Range("E16").Value = "  -1.26%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.33"
$ws.Range("E17").Value = "  +2.36%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.141.90"
$ws.Range("E18").Value = "  +0.38%  "
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "63.740.43"
$ws.Range("E19").Value = "  -0.13%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "469.47"
$ws.Range("E20").Value = "  +0.67%  "
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.732"
$ws.Range("E22").Value = "  -0.24%  "
$ws.Range("E23").Value = "  -0.60%  "
$ws.Range("E24").Value = "  +7.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.97"
$ws.Range("E25").Value = "  -2.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "81.47"
$ws.Range("E26").Value = "  -0.99%  "
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.87"
$ws.Range("E28").Value = "  +9.63%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.45"
$ws.Range("E29").Value = "  +8.47%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.71"
$ws.Range("E30").Value = "  -0.32%  "
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.24"
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("E32").Value = "  +0.18%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.70"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.110"
$ws.Range("E34").Value = "  +0.40%  "
$ws.Range("D35").Value = "0.0₃0841"
$ws.Range("E35").Value = "  -5.32%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.06"
$ws.Range("E36").Value = "  +0.74%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.31"
$ws.Range("E37").Value = "  -3.77%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.14"
$ws.Range("E38").Value = "  +0.48%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.21"
$ws.Range("E39").Value = "  -6.67%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "51.55"
$ws.Range("E40").Value = "  +1.19%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.32"
$ws.Range("E41").Value = "  +7.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "454.63"
$ws.Range("E42").Value = "  -0.47%  "
$ws.Range("E43").Value = "  +5.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0371"
$ws.Range("E44").Value = "  -0.63%  "
$ws.Range("D45").Value = "2.912.73"
$ws.Range("E45").Value = "  +0.66%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.26"
$ws.Range("E46").Value = "  +9.77%  "
$ws.Range("E47").Value = "  -3.20%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "130.84"
$ws.Range("E48").Value = "  +3.82%  "
$ws.Range("E49").Value = "  -0.04%  "
$ws.Range("E50").Value = "  +2.52%  "
$ws.Range("E51").Value = "  -1.23%  "
